$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 (B3, C3, D3)
$ws.Range("B3").Value = "長安路三段１６號"
$ws.Range("C3").Value = "查無結果"
$ws.Range("D3").Value = "查詢失敗"

# Update row 4 (B4) and row height
$ws.Range("B4").Value = "長安路016號"
$ws.Rows.Item(4).RowHeight = 15.75

# Delete rows 5-18
$ws.Range("A5:A18").EntireRow.Delete()

# Select A4:XFD4 (entire row 4)
$ws.Rows.Item(4).Select()
